$d = $word.ActiveDocument
$normalStyle = $d.Styles.Item("Normal")

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its current location (inside
#    the "a. Cat and seed?" paragraph). It gets re-created later at
#    the very end of the document, after the new content is typed.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the paragraph that contains the "B. The constraints..."
#    sentence - the three new paragraphs of possible solutions get
#    appended right after it.
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*The constraints are that some combination of items will result in one pet eating the other or the seed.*") {
        $targetIndex = $i
    }
}

$curIndex = $targetIndex

# Each entry is one new paragraph, expressed as an array of run texts
# that will be typed as individual runs.
$newParagraphs = @(
    @("       3. A.", " ", "Cat and seed- seems like correct choice"),
    @("                ", "Cat", (" and bird- leave cat or bird in cage so cat doesn" + [char]0x2019 + "t eat bird")),
    @("                ", "Bird", (" and seed- put bird in cage so bird doesn" + [char]0x2019 + "t eat seed."))
)

foreach ($runTexts in $newParagraphs) {
    $anchor = $d.Paragraphs.Item($curIndex).Range
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter() | Out-Null
    $curIndex = $curIndex + 1

    $firstRunParaIndex = $curIndex

    # Strip the inherited "ListParagraph" numbering/style immediately,
    # while the paragraph is still a single empty one - this keeps the
    # final merged paragraph free of any <w:pPr> element.
    $d.Paragraphs.Item($firstRunParaIndex).Range.ParagraphFormat.Style = $normalStyle

    $segIndex = 0
    foreach ($segText in $runTexts) {
        $segRange = $d.Paragraphs.Item($curIndex).Range
        $segRange.Collapse(0)
        $segRange.InsertAfter($segText)

        if ($segIndex -lt ($runTexts.Count - 1)) {
            # Open a fresh paragraph for the next run so that it cannot
            # merge with the one we just typed, then splice the paragraph
            # mark back out - this keeps the two chunks of text as
            # separate <w:r> runs once the mark is gone.
            $segRange2 = $d.Paragraphs.Item($curIndex).Range
            $segRange2.Collapse(0)
            $segRange2.InsertParagraphAfter() | Out-Null
            $curIndex = $curIndex + 1
        }
        $segIndex = $segIndex + 1
    }

    # Merge all the single-run paragraphs we just created (from
    # $firstRunParaIndex to $curIndex) back into one paragraph by
    # deleting the paragraph marks between them.
    while ($curIndex -gt $firstRunParaIndex) {
        $prevPara = $d.Paragraphs.Item($firstRunParaIndex)
        $markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
        $markRange.Delete()
        $curIndex = $curIndex - 1
    }

    $curIndex = $firstRunParaIndex
}

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark, collapsed, at the very end of
#    the document (right after the last run of new text, before the
#    final paragraph mark). A temporary marker character is used to
#    sidestep an engine quirk where a bookmark collapsed exactly on a
#    paragraph-end boundary snaps back to the start of the document.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter([char]0x2060)   # temporary zero-width marker char

$markerPos = $d.Content.End - 2
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$markerRange.Delete()
